# Modifications pour utiliser XGBClassifier et ajuster les predictions
# - Renames the S+1/S+2/S+3 headers on "Valeurs reelles" to add a "_class" suffix
# - Replaces the continuous price forecasts with small integer class labels
#   (XGBClassifier outputs) on both the "Valeurs reelles" and "Predictions" sheets

$wb = $excel.ActiveWorkbook
$wsReelles = $wb.Worksheets.Item("Valeurs réelles")
$wsPred = $wb.Worksheets.Item("Prédictions")

# --- Header renames on "Valeurs reelles" (append "_class") ---
$wsReelles.Range("C1").Value = "PRIX EXP POMME GOLDEN FRANCE 170/220G CAT.I PLATEAU 1RG_S+1_class"
$wsReelles.Range("D1").Value = "PRIX EXP POMME GOLDEN FRANCE 170/220G CAT.I PLATEAU 1RG_S+2_class"
$wsReelles.Range("E1").Value = "PRIX EXP POMME GOLDEN FRANCE 170/220G CAT.I PLATEAU 1RG_S+3_class"

# --- New classifier-style integer values for columns C:E, rows 2-28 ---
$reellesValues = @{
    "C2" = 4
    "D2" = 2
    "E2" = 2
    "C3" = 2
    "D3" = 2
    "E3" = 2
    "C4" = 2
    "D4" = 2
    "E4" = 2
    "C5" = 2
    "D5" = 2
    "E5" = 2
    "C6" = 2
    "D6" = 2
    "E6" = 2
    "C7" = 2
    "D7" = 2
    "E7" = 2
    "C8" = 2
    "D8" = 2
    "E8" = 2
    "C9" = 2
    "D9" = 2
    "E9" = 2
    "C10" = 2
    "D10" = 2
    "E10" = 1
    "C11" = 2
    "D11" = 1
    "E11" = 1
    "C12" = 1
    "D12" = 1
    "E12" = 2
    "C13" = 1
    "D13" = 2
    "E13" = 2
    "C14" = 2
    "D14" = 2
    "E14" = 0
    "C15" = 2
    "D15" = 0
    "E15" = 3
    "C16" = 0
    "D16" = 3
    "E16" = 1
    "C17" = 3
    "D17" = 1
    "E17" = 1
    "C18" = 1
    "D18" = 1
    "E18" = 2
    "C19" = 1
    "D19" = 2
    "E19" = 2
    "C20" = 2
    "D20" = 2
    "E20" = 2
    "C21" = 2
    "D21" = 2
    "E21" = 3
    "C22" = 2
    "D22" = 3
    "E22" = 2
    "C23" = 3
    "D23" = 2
    "E23" = 2
    "C24" = 2
    "D24" = 2
    "E24" = 2
    "C25" = 2
    "D25" = 2
    "E25" = 2
    "C26" = 2
    "D26" = 2
    "E26" = 2
    "C27" = 2
    "D27" = 2
    "E27" = 2
    "C28" = 2
    "D28" = 2
    "E28" = 2
}

foreach ($addr in $reellesValues.Keys) {
    $wsReelles.Range($addr).Value = $reellesValues[$addr]
}

# --- New classifier-style integer values for columns B:D, rows 2-28 on "Predictions" ---
$predValues = @{
    "B2" = 0
    "C2" = 0
    "D2" = 1
    "B3" = 0
    "C3" = 2
    "D3" = 0
    "B4" = 0
    "C4" = 1
    "D4" = 1
    "B5" = 0
    "C5" = 2
    "D5" = 0
    "B6" = 0
    "C6" = 2
    "D6" = 2
    "B7" = 0
    "C7" = 2
    "D7" = 0
    "B8" = 2
    "C8" = -1
    "D8" = 2
    "B9" = -2
    "C9" = 0
    "D9" = 2
    "B10" = 2
    "C10" = 2
    "D10" = 2
    "B11" = -2
    "C11" = -1
    "D11" = 2
    "B12" = 2
    "C12" = 0
    "D12" = 2
    "B13" = -2
    "C13" = -1
    "D13" = 0
    "B14" = 1
    "C14" = 0
    "D14" = 0
    "B15" = 0
    "C15" = 0
    "D15" = 0
    "B16" = 0
    "C16" = 0
    "D16" = 0
    "B17" = -1
    "C17" = 0
    "D17" = 0
    "B18" = 0
    "C18" = 0
    "D18" = 0
    "B19" = 0
    "C19" = 0
    "D19" = 0
    "B20" = 0
    "C20" = 0
    "D20" = -1
    "B21" = 0
    "C21" = 0
    "D21" = -1
    "B22" = 0
    "C22" = 0
    "D22" = -1
    "B23" = 2
    "C23" = 0
    "D23" = 0
    "B24" = 1
    "C24" = 0
    "D24" = 0
    "B25" = 1
    "C25" = -1
    "D25" = 0
    "B26" = 0
    "C26" = 0
    "D26" = 2
    "B27" = 0
    "C27" = -1
    "D27" = 0
    "B28" = 0
    "C28" = 0
    "D28" = 0
}

foreach ($addr in $predValues.Keys) {
    $wsPred.Range($addr).Value = $predValues[$addr]
}
